$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "23-03-2025"
$ws.Range("B4").Value = "Sunrisers Hyderabad vs Rajasthan Royals"
$ws.Range("C4").Value = "Sunrisers Hyderabad"
$ws.Range("D4").Value = "Sunrisers Hyderabad"
$ws.Range("E4").Value = "nandini"

$ws.Range("A5").Value = "23-03-2025"
$ws.Range("B5").Value = "Chennai Super Kings vs Mumbai Indians"
$ws.Range("C5").Value = "Chennai Super Kings"
$ws.Range("D5").Value = "Chennai Super Kings"
$ws.Range("E5").Value = "nandini"
